# Update cryptocurrency price/volume table cells per the Wed Jul 26 21:25:46 UTC 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.548.70"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "1.878.12"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'0.7246"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "'239.66"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07876"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").Value = "'0.3093"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "'25.25"
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("D11").Value = "'0.08224"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7268"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.264"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.849.50"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'90.42"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "29.632.09"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "'5.867"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "'0.000007881"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'242.61"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "'13.38"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "2.142.10"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'0.9991"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'7.778"
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("D25").Value = "'0.1598"
$ws.Range("E25").Value = "  +11.18%  "
$ws.Range("D26").Value = "'162.85"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'8.998"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "'18.40"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'1.948"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "'1.368"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "'1.482"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'4.357"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'4.094"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "'0.05269"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").Value = "'1.198"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").Value = "'0.7186"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'0.01871"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "'2.717"
$ws.Range("D41").Value = "1.188.72"
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("D42").Value = "'0.8999"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "'6.013"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4327"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'71.79"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").Value = "'0.9989"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'102.82"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'0.5359"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "'1.779"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'9.250"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.089"
$ws.Range("E51").Value = "  +1.91%  "
